# Belgium First Division A - base update (30-05-2024)
#
# The underlying source data was re-sorted, which moved some match rows
# relative to each other. For the affected groups of rows, the entire
# row content (columns B:AD -- everything except the running "id" in
# column A) needs to be rotated among the rows in the group so that each
# row ends up holding the data that used to belong to the "next" row in
# the group (wrapping around for groups with more than two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array lists the row numbers (in sheet order) that form one
# rotation group. Row i receives the B:AD values that row (i+1) held
# before the edit; the last row in a group wraps around to the first.
$groups = @(
    ,@(2, 3)
    ,@(26, 27)
    ,@(155, 156)
    ,@(185, 186)
    ,@(241, 242, 243, 244)
    ,@(276, 277)
    ,@(297, 298)
)

foreach ($group in $groups) {
    $n = $group.Length

    # Snapshot the current B:AD values for every row in this group first,
    # so later writes in this loop don't clobber data we still need to read.
    $snapshots = @()
    foreach ($row in $group) {
        $rng = $ws.Range("B$row`:AD$row")
        $snapshots += ,$rng.Value2
    }

    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $group[$i]
        $srcIndex = ($i + 1) % $n
        $destRng = $ws.Range("B$destRow`:AD$destRow")
        $destRng.Value2 = $snapshots[$srcIndex]
    }
}
